$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date in column C was bumped from 45171 to 45172
# for every data row (rows 2 through 54).
$ws.Range("C2:C54").Value = 45172
